$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.208.07"
$ws.Range("E2").Value = "  +4.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.650.09"
$ws.Range("E3").Value = "  +11.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.04"
$ws.Range("E5").Value = "  +7.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.52"
$ws.Range("E6").Value = "  +11.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").Value = "  +10.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +20.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.97"
$ws.Range("E10").Value = "  +17.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.09"
$ws.Range("E11").Value = "  +4.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0859"
$ws.Range("E12").Value = "  +10.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.51"
$ws.Range("E13").Value = "  +22.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.055.89"
$ws.Range("E14").Value = "  +11.49%  "
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.653.38"
$ws.Range("E16").Value = "  +11.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.952"
$ws.Range("E17").Value = "  +15.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.54"
$ws.Range("E18").Value = "  +11.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "47.952.33"
$ws.Range("E19").Value = "  +6.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").Value = "  +12.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("E21").Value = "  +8.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.90"
$ws.Range("E22").Value = "  +13.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.10"
$ws.Range("E23").Value = "  +10.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "273.29"
$ws.Range("E24").Value = "  +14.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  +14.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "31.22"
$ws.Range("E26").Value = "  +49.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.22"
$ws.Range("E27").Value = "  +18.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.07"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.80"
$ws.Range("E30").Value = "  +13.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "40.97"
$ws.Range("E31").Value = "  +7.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.31"
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.30"
$ws.Range("E33").Value = "  +17.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.80"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.31"
$ws.Range("E35").Value = "  +19.63%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.90"
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0861"
$ws.Range("E37").Value = "  +13.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "153.75"
$ws.Range("E38").Value = "  +4.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.123"
$ws.Range("E39").Value = "  +10.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("E40").Value = "  +9.59%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.28"
$ws.Range("E41").Value = "  +61.53%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.02"
$ws.Range("E42").Value = "  +15.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.35"
$ws.Range("E43").Value = "  +16.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.74"
$ws.Range("E44").Value = "  +18.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0335"
$ws.Range("E45").Value = "  +14.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.192.29"
$ws.Range("E46").Value = "  +12.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.79"
$ws.Range("E47").Value = "  +9.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.99"
$ws.Range("E49").Value = "  +17.49%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.06"
$ws.Range("E50").Value = "  +16.94%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.87"
$ws.Range("E51").Value = "  +8.69%  "
